$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = -7.263000000000001
$ws.Range("A12").Value = -21.542
$ws.Range("D23").Value = -8.082000000000001
$ws.Range("D28").Value = -8.406000000000001
$ws.Range("A32").Value = -21.591
$ws.Range("D32").Value = -7.306999999999999
$ws.Range("D34").Value = -7.933
$ws.Range("A36").Value = -20.316
$ws.Range("A38").Value = -19.83
$ws.Range("D42").Value = -8.350999999999999
$ws.Range("A46").Value = -21.823
$ws.Range("A54").Value = -22.141
$ws.Range("D54").Value = -7.877000000000001
$ws.Range("A55").Value = -22.184
$ws.Range("A67").Value = -21.536
$ws.Range("A69").Value = -21.503
$ws.Range("A72").Value = -21.689
$ws.Range("A91").Value = -20.675
$ws.Range("D97").Value = -7.674000000000001
$ws.Range("A99").Value = -22.008
$ws.Range("D99").Value = -8.280000000000001
$ws.Range("D101").Value = -7.784999999999999
$ws.Range("A104").Value = -21.385
